# g2 gas pressure added
#
# 1. Bump the "Date" timestamp in the header table.
# 2. In the "DIRECT FIRED CIRCUIT" table: swap the last-column values of rows
#    "2." (Fuel Type) and "3." (Calorific value type), then insert a new row
#    after "3." for "Gas Pressure" / "mbar" / "100".
# 3. Reword the final "Burner Selection" note.

$d = $word.ActiveDocument

# 1. Update the enquiry date/time stamp.
$d.Content.Find.Execute("30-Jun-2021, 16:00 ", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "08-Jul-2021, 15:29 ", 2)

# 2. Locate the "DIRECT FIRED CIRCUIT" table (5th table in the document).
$tbl = $d.Tables.Item(5)

# Row 3 ("2.  Fuel Type") last cell: Normal -> NaturalGas
$tbl.Rows.Item(3).Cells.Item(4).Range.Text = "NaturalGas"
# Row 4 ("3.  Calorific value type") last cell: NaturalGas -> Normal
$tbl.Rows.Item(4).Cells.Item(4).Range.Text = "Normal"

# Insert a new row ("Gas Pressure") right before row "4. Calorific Value".
$calorificValueRow = $tbl.Rows.Item(5)
$newRow = $tbl.Rows.Add($calorificValueRow)
$newRow.Cells.Item(2).Range.Text = "Gas Pressure"
$newRow.Cells.Item(3).Range.Text = "mbar"
$newRow.Cells.Item(4).Range.Text = "100"

# 3. Update the closing note about Burner Selection.
$d.Content.Find.Execute("5. Burner Selection is for Sea level", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "5. Burner Selection is valid upto 100m above mean Sea level.", 2)
